$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use row 71 (a fully-populated row with the standard styles used across the
# History sheet) as the style template for the two new rows being appended.
$styleRow = 71

# --- Row 73 -----------------------------------------------------------
$ws.Cells.Item($styleRow, 1).Copy()
$ws.Cells.Item(73, 1).PasteSpecial(-4122)
$ws.Cells.Item(73, 1).Value = 45304

$ws.Cells.Item($styleRow, 2).Copy()
$ws.Cells.Item(73, 2).PasteSpecial(-4122)
$ws.Cells.Item(73, 2).Value = 0.1491550925925926

$ws.Cells.Item($styleRow, 3).Copy()
$ws.Cells.Item(73, 3).PasteSpecial(-4122)
$ws.Cells.Item(73, 3).Value = "Программирование"

$ws.Cells.Item($styleRow, 4).Copy()
$ws.Cells.Item(73, 4).PasteSpecial(-4122)
$ws.Cells.Item(73, 4).Value = "Добавил проверку дня недели если текущее время меньше времени последней ячейки"

# --- Row 74 -----------------------------------------------------------
$ws.Cells.Item($styleRow, 1).Copy()
$ws.Cells.Item(74, 1).PasteSpecial(-4122)
$ws.Cells.Item(74, 1).Value = 45304

$ws.Cells.Item($styleRow, 2).Copy()
$ws.Cells.Item(74, 2).PasteSpecial(-4122)
$ws.Cells.Item(74, 2).Value = 0.18293981481481481

$ws.Cells.Item($styleRow, 3).Copy()
$ws.Cells.Item(74, 3).PasteSpecial(-4122)
$ws.Cells.Item(74, 3).Value = "Программирование"

# Widen column D to fit the new, longer history note. (Excel stores/serialises
# column widths snapped to whole pixels, so the ColumnWidth value below is the
# one that round-trips to the closest possible width to the target 101.77734375.)
$ws.Columns.Item(4).ColumnWidth = 100.94401041666667

# Update the view so the newly added rows are visible/selected, mirroring
# what Excel records after a user scrolls to and selects the new cell.
$ws.Application.ActiveWindow.ScrollRow = 45
$ws.Range("D74").Select()
